$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 251, pushing the existing data (old rows 251-306)
# down to rows 255-310.
$ws.Rows("251:254").Insert()

# Constant values shared by every data row in this sheet (columns A,B,C,E-K).
$mercadoId = 6
$mercado = "Mercado Mayorista Lo Valledor de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100107
$producto = "Otros"
$categoriaId = 100107011
$categoria = "Tuna"
$variedad = "Sin especificar"

# Data for the 4 newly inserted rows (251-254).
$newRows = @(
    @{ Row=251; Fecha=44641; Calidad="Especial";                 Volumen=150; PMin=13000; PMax=13000; PProm=13000; Unidad="$/caja 18 kilos"; Origen="Provincia de Melipilla"; PKg=722; KgUnidad=18 },
    @{ Row=252; Fecha=44641; Calidad="Extra (doble especial)";   Volumen=140; PMin=15000; PMax=15000; PProm=15000; Unidad="$/caja 18 kilos"; Origen="Provincia de Melipilla"; PKg=833; KgUnidad=18 },
    @{ Row=253; Fecha=44641; Calidad="Primera";                  Volumen=170; PMin=12000; PMax=12000; PProm=12000; Unidad="$/caja 18 kilos"; Origen="Provincia de Melipilla"; PKg=667; KgUnidad=18 },
    @{ Row=254; Fecha=44641; Calidad="Segunda";                  Volumen=140; PMin=8000;  PMax=8000;  PProm=8000;  Unidad="$/caja 18 kilos"; Origen="Provincia de Melipilla"; PKg=444; KgUnidad=18 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.PMin
    $ws.Cells.Item($row, 15).Value = $r.PMax
    $ws.Cells.Item($row, 16).Value = $r.PProm
    $ws.Cells.Item($row, 17).Value = $r.Unidad
    $ws.Cells.Item($row, 18).Value = $r.Origen
    $ws.Cells.Item($row, 19).Value = $r.PKg
    $ws.Cells.Item($row, 20).Value = $r.KgUnidad
}
